$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

$ws.Rows.Item(87).Insert() | Out-Null
$ws.Rows.Item(99).Cut($ws.Rows.Item(87)) | Out-Null
$ws.Rows.Item(99).Delete() | Out-Null

Write-Host "After insert+cut+delete:"
for ($r = 85; $r -le 101; $r++) {
    Write-Host "Row $r A=$($ws.Range("A$r").Value()) B=$($ws.Range("B$r").Value()) height=$($ws.Rows.Item($r).RowHeight)"
}
$ur = $ws.UsedRange.Address()
Write-Host "UsedRange: $ur"
for ($r = 268; $r -le 272; $r++) {
    Write-Host "Row $r A=$($ws.Range("A$r").Value())"
}
